# RR_Actualización_productos_de_información.xlsx
# "se ajusta para poder exportar excel con productos de infromacion"
#
# - Add two new header columns (Enlace rt, Tipo de emergencia) and rename
#   the existing "Enlace" header to "Enlace rr"
# - Give the new/renamed link headers (H1:J1) the same header formatting
#   used by the rest of row 1
# - Remove the two example hyperlinks that had been left in row 2 (H2/I2)
#   together with the "por favor respetar este formato" note in J2
# - Move the active selection back to A2

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing header formatting onto the new header cells so H1:J1
# match the look of the rest of the header row (bold white text on the
# colored fill, centered + wrapped).
$ws.Range("G1").Copy()
$ws.Range("H1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Header row text
$ws.Range("A1").Value = "Tipo de Respuesta"
$ws.Range("B1").Value = "Tipo de producto"
$ws.Range("F1").Value = "Fecha de Elaboración"
$ws.Range("H1").Value = "Enlace rr"
$ws.Range("I1").Value = "Enlace rt"
$ws.Range("J1").Value = "Tipo de emergencia"

# Drop the leftover example hyperlinks and their backing text.
$ws.Hyperlinks.Delete()
$ws.Range("H2").ClearContents()
$ws.Range("I2").ClearContents()
$ws.Range("J2").ClearContents()

# Reset the selection to A2.
$ws.Range("A2").Select()
